$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

$para = $tr.Paragraphs(3)
$run = $para.Runs(1)
$run.Text = "Sprechstunde: Dienstag 12-14 Uhr "
